$wb = $excel.ActiveWorkbook

# --- Sheet "Journal de travail HAY": log a new entry of work done ---
$journalHay = $wb.Worksheets.Item("Journal de travail HAY")
$null = $journalHay.Activate()
$journalHay.Range("A11").Value = 45022
$journalHay.Range("B11").Value = "REST des comptes"
$journalHay.Range("C11").Value = 1.5
$null = $journalHay.Range("B17").Select()

# --- Sheet "Journal de travail CLA": log a new entry of work done ---
$journalCla = $wb.Worksheets.Item("Journal de travail CLA")
$null = $journalCla.Activate()
$journalCla.Range("A11").Value = 45022
$journalCla.Range("B11").Value = "REST des livres"
$journalCla.Range("C11").Value = 1.5
$null = $journalCla.Range("C13").Select()

# --- Sheet "Planning": mark column H (Bloc 3 / 06.03.2023 - PM) status for a few tasks ---
$planning = $wb.Worksheets.Item("Planning")
$null = $planning.Activate()
$planning.Range("H18").Value = "A"
$planning.Range("H19").Value = "N"
$planning.Range("H20").Value = "N"

# Update the visible selection to match the latest working view (keeps Planning the active tab)
$null = $planning.Range("G25").Select()

$excel.Calculate()
